# The commit swaps the presentation's colour theme from the custom
# "Integral" scheme to the stock PowerPoint "Office" scheme (the notes
# master keeps a theme part too, but that part is not something the
# PowerPoint object model lets an end-user / macro touch directly -
# only the colours that are actually applied to the deck, which live on
# the slide master's theme, are reachable here).
#
# PowerPoint COM exposes the 12 theme colour slots through
# Master.ColorScheme.Colors(index).RGB, where index 1-12 maps to
# dk1, lt1, dk2, lt2, accent1..accent6, hlink, folHlink (in that
# order) and RGB values are packed as 0x00BBGGRR (the standard VBA
# RGB() colour order).
#
# Target values come from the stock Office theme colour scheme:
#   dk1      #000000
#   lt1      #FFFFFF
#   dk2      #44546A
#   lt2      #E7E6E6
#   accent1  #5B9BD5
#   accent2  #ED7D31
#   accent3  #A5A5A5
#   accent4  #FFC000
#   accent5  #4472C4
#   accent6  #70AD47
#   hlink    #0563C1
#   folHlink #954F72

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.ColorScheme

$cs.Colors(1).RGB  = 0          # dk1      #000000
$cs.Colors(2).RGB  = 16777215   # lt1      #FFFFFF
$cs.Colors(3).RGB  = 6968388    # dk2      #44546A
$cs.Colors(4).RGB  = 15132391   # lt2      #E7E6E6
$cs.Colors(5).RGB  = 13998939   # accent1  #5B9BD5
$cs.Colors(6).RGB  = 3243501    # accent2  #ED7D31
$cs.Colors(7).RGB  = 10855845   # accent3  #A5A5A5
$cs.Colors(8).RGB  = 49407      # accent4  #FFC000
$cs.Colors(9).RGB  = 12874308   # accent5  #4472C4
$cs.Colors(10).RGB = 4697456    # accent6  #70AD47
$cs.Colors(11).RGB = 12673797   # hlink    #0563C1
$cs.Colors(12).RGB = 7491477    # folHlink #954F72
